# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-31 01:24:38
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet contains
# comma-separated lists of recorder names/emails. This edit rotates the first
# entry of each such list to the end (i.e. "A, B" -> "B, A") for the specific
# combinations that were re-ordered upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact (old -> new) text replacements observed for column G ("Recorded By").
# Using a literal lookup (rather than a generic "always rotate" rule) because
# only cells with these precise values were touched upstream; other multi-value
# combinations (e.g. "backup@backdoor.com, System", "System, admin@admin.com")
# were left unchanged.
$map = @{
    "system, backup@backdoor.com, System" = "backup@backdoor.com, System, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
